$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "G2" = 22.628972
    "H2" = 67.886916
    "I2" = 0.004372730881336598
    "J2" = 0.004372730881336598
    "M2" = 3.319447
    "N2" = 9.958341000000001
    "O2" = 0.03276237985858125
    "P2" = 0.03276237985858125
    "Q2" = 75.11567321848401
    "R2" = 676.041058966356
    "S2" = 0.0001432610701536984
    "T2" = 0.0001432610701536984
    "G3" = 22.628972
    "H3" = 67.886916
    "I3" = 0.004372730881336598
    "J3" = 0.004372730881336598
    "O3" = 0.7010975337260504
    "P3" = 0.7010975337260504
    "Q3" = 1607.435524066707
    "R3" = 14466.91971660036
    "S3" = 0.003065710836552828
    "T3" = 0.003065710836552828
    "G4" = 22.628972
    "H4" = 67.886916
    "I4" = 0.004372730881336598
    "J4" = 0.004372730881336598
    "M4" = 26.57769466666667
    "N4" = 79.73308400000001
    "O4" = 0.262317346363633
    "P4" = 0.262317346363633
    "Q4" = 601.4259084365494
    "R4" = 5412.833175928944
    "S4" = 0.001147043161154527
    "T4" = 0.001147043161154527
    "G5" = 22.628972
    "H5" = 67.886916
    "I5" = 0.004372730881336598
    "J5" = 0.004372730881336598
    "M5" = 0.3873156666666667
    "N5" = 1.161947
    "O5" = 0.003822740051735415
    "P5" = 0.003822740051735415
    "Q5" = 8.764555376161335
    "R5" = 78.880998385452
    "S5" = 0.00001671581347554572
    "T5" = 0.00001671581347554572
    "I6" = 0.9598063873258337
    "J6" = 0.9598063873258338
    "M6" = 3.319447
    "N6" = 9.958341000000001
    "O6" = 0.03276237985858125
    "P6" = 0.03276237985858125
    "Q6" = 16487.75213930923
    "R6" = 148389.769253783
    "S6" = 0.03144554145226153
    "T6" = 0.03144554145226153
    "I7" = 0.9598063873258337
    "J7" = 0.9598063873258338
    "O7" = 0.7010975337260504
    "P7" = 0.7010975337260504
    "S7" = 0.6729178910086523
    "T7" = 0.6729178910086524
    "I8" = 0.9598063873258337
    "J8" = 0.9598063873258338
    "M8" = 26.57769466666667
    "N8" = 79.73308400000001
    "O8" = 0.262317346363633
    "P8" = 0.262317346363633
    "Q8" = 132011.8809242144
    "R8" = 1188106.92831793
    "S8" = 0.251773864546178
    "T8" = 0.251773864546178
    "I9" = 0.9598063873258337
    "J9" = 0.9598063873258338
    "M9" = 0.3873156666666667
    "N9" = 1.161947
    "O9" = 0.003822740051735415
    "P9" = 0.003822740051735415
    "Q9" = 1923.803787700576
    "R9" = 17314.23408930518
    "S9" = 0.003669090318741939
    "T9" = 0.00366909031874194
    "G10" = 182.6322073333333
    "H10" = 547.896622
    "I10" = 0.0352911079183418
    "J10" = 0.0352911079183418
    "M10" = 3.319447
    "N10" = 9.958341000000001
    "O10" = 0.03276237985858125
    "P10" = 0.03276237985858125
    "Q10" = 606.2379327360113
    "R10" = 5456.141394624102
    "S10" = 0.001156220683250899
    "T10" = 0.001156220683250899
    "G11" = 182.6322073333333
    "H11" = 547.896622
    "I11" = 0.0352911079183418
    "J11" = 0.0352911079183418
    "O11" = 0.7010975337260504
    "P11" = 0.7010975337260504
    "Q11" = 12973.16987737296
    "R11" = 116758.5288963566
    "S11" = 0.02474250872400932
    "T11" = 0.02474250872400932
    "G12" = 182.6322073333333
    "H12" = 547.896622
    "I12" = 0.0352911079183418
    "J12" = 0.0352911079183418
    "M12" = 26.57769466666667
    "N12" = 79.73308400000001
    "O12" = 0.262317346363633
    "P12" = 0.262317346363633
    "Q12" = 4853.943042804694
    "R12" = 43685.48738524225
    "S12" = 0.009257469779372016
    "T12" = 0.009257469779372016
    "G13" = 182.6322073333333
    "H13" = 547.896622
    "I13" = 0.0352911079183418
    "J13" = 0.0352911079183418
    "M13" = 0.3873156666666667
    "N13" = 1.161947
    "O13" = 0.003822740051735415
    "P13" = 0.003822740051735415
    "Q13" = 70.73631513811489
    "R13" = 636.6268362430341
    "S13" = 0.000134908731709562
    "T13" = 0.000134908731709562
    "G14" = 2.741590666666667
    "H14" = 8.224772
    "I14" = 0.00052977387448787
    "J14" = 0.00052977387448787
    "M14" = 3.319447
    "N14" = 9.958341000000001
    "O14" = 0.03276237985858125
    "P14" = 0.03276237985858125
    "Q14" = 9.100564913694667
    "R14" = 81.905084223252
    "S14" = 0.00001735665291512394
    "T14" = 0.00001735665291512394
    "G15" = 2.741590666666667
    "H15" = 8.224772
    "I15" = 0.00052977387448787
    "J15" = 0.00052977387448787
    "O15" = 0.7010975337260504
    "P15" = 0.7010975337260504
    "Q15" = 194.7472571909022
    "R15" = 1752.72531471812
    "S15" = 0.0003714231568359399
    "T15" = 0.0003714231568359399
    "G16" = 2.741590666666667
    "H16" = 8.224772
    "I16" = 0.00052977387448787
    "J16" = 0.00052977387448787
    "M16" = 26.57769466666667
    "N16" = 79.73308400000001
    "O16" = 0.262317346363633
    "P16" = 0.262317346363633
    "Q16" = 72.86515963964979
    "R16" = 655.786436756848
    "S16" = 0.0001389688769284384
    "T16" = 0.0001389688769284384
    "G17" = 2.741590666666667
    "H17" = 8.224772
    "I17" = 0.00052977387448787
    "J17" = 0.00052977387448787
    "M17" = 0.3873156666666667
    "N17" = 1.161947
    "O17" = 0.003822740051735415
    "P17" = 0.003822740051735415
    "Q17" = 1.061861016787111
    "R17" = 9.556749151084
    "S17" = 0.000002025187808367831
    "T17" = 0.000002025187808367832
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}

Write-Host "Updated $($newValues.Count) cells"
